$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the price/volume columns so purely numeric-looking
# strings (e.g. '299.80') are not auto-converted to numbers by Excel's
# type inference, matching the inlineStr text cells in the source file.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '23.457.93'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.644.67'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '299.80'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '0.3785'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('D8').Value = '50.52'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('D10').Value = '0.08050'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '22.07'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').Value = '6.288'
$ws.Range('E14').Value = '  -3.13%  '
$ws.Range('D15').Value = '7.239'
$ws.Range('E15').Value = '  -3.61%  '
$ws.Range('D16').Value = '0.00001207'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').Value = '1.648.99'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').Value = '95.21'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').Value = '0.06977'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '6.615'
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '23.469.03'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = '2.417'
$ws.Range('E25').Value = '  -4.21%  '
$ws.Range('D26').Value = '3.004'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').Value = '21.04'
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = '151.81'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = '131.47'
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('D31').Value = '1.828.92'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = '6.848'
$ws.Range('E32').Value = '  -5.24%  '
$ws.Range('E33').Value = '  -5.02%  '
$ws.Range('E34').Value = '  -7.63%  '
$ws.Range('D35').Value = '0.9896'
$ws.Range('E35').Value = '  -6.54%  '
$ws.Range('D36').Value = '0.02688'
$ws.Range('E36').Value = '  -4.36%  '
$ws.Range('D37').Value = '0.08739'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').Value = '5.906'
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('D39').Value = '0.2418'
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').Value = '0.06790'
$ws.Range('E40').Value = '  -3.23%  '
$ws.Range('D41').Value = '12.83'
$ws.Range('E41').Value = '  -3.17%  '
$ws.Range('D42').Value = '0.6886'
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('D44').Value = '15.55'
$ws.Range('E44').Value = '  -2.79%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '0.6388'
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('D47').Value = '3.921'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').Value = '2.240'
$ws.Range('E48').Value = '  -2.99%  '
$ws.Range('D49').Value = '127.25'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').Value = '0.07668'
$ws.Range('E50').Value = '  -3.26%  '
$ws.Range('D51').Value = '1.239'
$ws.Range('E51').Value = '  +3.00%  '

# Restore the default (Normal) cell style so no stray number format is
# left behind on the cells, matching the original formatting.
$dataRange.Style = "Normal"
